# Renames the header row from "<name>_old"/"<name>_new" to
# "<name>_FV2210"/"<name>_FV2304" respectively (column K, "diff", is left
# untouched), turns the sheet's used range into a real Excel Table
# ("Table1") so the new header names become the table's column names, and
# freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (A1:U1) ------------------------------------------
$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)

$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

# Columns A..J (1..10)
for ($i = 0; $i -lt $headersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2210[$i]
}

# Column K (11) stays "diff" - untouched on purpose.

# Columns L..U (12..21)
for ($i = 0; $i -lt $headersFV2304.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2304[$i]
}

# --- 2. Convert the data range into an Excel Table (Table1) ----------------
$tableRange = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $tableRange, 0, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null

Write-Host "Renamed header row, added Table1 (A1:U79) and froze the top row."
